# Update the division-problem answers in the first (and only) table.
# Each data row contains 5 cells; data rows are 1, 5, 9, 13, 17 (1-indexed).
# The replacement values are applied in row-major, left-to-right order,
# matching the order of changes in the target diff.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "63÷8=7, 7", "87÷9=9, 6", "63÷8=7, 7", "89÷9=9, 8", "29÷2=14, 1",
    "23÷3=7, 2", "12÷5=2, 2", "12÷7=1, 5", "98÷2=49, 0", "99÷2=49, 1",
    "92÷9=10, 2", "17÷5=3, 2", "37÷5=7, 2", "64÷3=21, 1", "25÷7=3, 4",
    "91÷3=30, 1", "30÷6=5, 0", "67÷2=33, 1", "35÷8=4, 3", "42÷9=4, 6",
    "68÷9=7, 5", "88÷3=29, 1", "75÷6=12, 3", "40÷5=8, 0", "94÷9=10, 4"
)

$dataRows = @(1, 5, 9, 13, 17)

$idx = 0
foreach ($r in $dataRows) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated $idx cells"
